$d = $word.ActiveDocument

# The title paragraph (paragraph 1) is the anchor: InsertXML on a collapsed
# range re-materializes the paragraph that contains the insertion point, so
# we re-supply its own content/formatting unchanged and append the three new
# bold paragraphs right after it. Everything after (the existing blank
# paragraph and the "No assumptions made" paragraph) is left untouched.
$r = $d.Paragraphs(1).Range
$r.Collapse(0)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="67543B9A" w14:textId="4E4E4F68" w:rsidR="006E1334" w:rsidRDefault="007D2A8C" w:rsidP="007D2A8C"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t>Assignment 1A Assumptions</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Bassel </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Emadeldin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Hamed </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Hamed</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Abdelkader</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>29796776</w:t></w:r></w:p>
'@

$null = $r.InsertXML($xml)
